# "completed E5 titrations 0329"
# Append the new titration result row (row 73) to CRMAccuracyData.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRMAccuracyData")

$newRow = 73

$ws.Cells.Item($newRow, 1).Value = 20220329
$ws.Cells.Item($newRow, 2).Value = 2234.105
$ws.Cells.Item($newRow, 3).Value = 2224.4699999999998
$ws.Range("D" + $newRow).Formula = "=100*(B" + $newRow + "-C" + $newRow + ")/C" + $newRow
$ws.Cells.Item($newRow, 5).Value = 180
$ws.Cells.Item($newRow, 6).Value = "CRM OPENED 20220318"

# Scroll the view down and move the selection the way the author left it
# (new bottom of the data, just past the freshly-entered row).
$win = $excel.ActiveWindow
$win.ScrollRow = 58
$win.ScrollColumn = 1
$ws.Range("D75").Select()

# Reposition the workbook window on screen, matching the saved workbookView.
$win.Left = 3750
